$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 data: last-run timestamp, position, and "found on first page" flag
$ws.Range("B3").Value = "06.02.2019, 22:42:39"
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = "'True"

# Move the active selection to G17
$ws.Range("G17").Select()
